$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- paises.xlsx update: refreshed COVID-19 country stats + provincias Spain snapshot (08:20) ---

# Row 1: refresh the "last updated" timestamp label
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Abril de 2020 a las 08:20"

# Row 16: updated case totals
$ws.Cells.Item(16, 2).Value = 10769
$ws.Cells.Item(16, 3).Value = 58
$ws.Cells.Item(16, 5).Value = 9187

# Row 37: updated case totals
$ws.Cells.Item(37, 2).Value = 2291
$ws.Cells.Item(37, 3).Value = 173
$ws.Cells.Item(37, 4).Value = 107
$ws.Cells.Item(37, 5).Value = 2153

# Row 67: updated case totals
$ws.Cells.Item(67, 2).Value = 649
$ws.Cells.Item(67, 3).Value = 68
$ws.Cells.Item(67, 5).Value = 634

# Row 68 (Armenia -> Hungria): updated case totals
$ws.Cells.Item(68, 1).Value = "Hungria"
$ws.Cells.Item(68, 2).Value = 585
$ws.Cells.Item(68, 3).Value = 60
$ws.Cells.Item(68, 4).Value = 42
$ws.Cells.Item(68, 5).Value = 522
$ws.Cells.Item(68, 6).Value = 17
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 21

# Row 69 (Barein -> Armenia): updated case totals
$ws.Cells.Item(69, 1).Value = "Armenia"
$ws.Cells.Item(69, 2).Value = 571
$ws.Cells.Item(69, 4).Value = 31
$ws.Cells.Item(69, 5).Value = 536
$ws.Cells.Item(69, 6).Value = 30

# Row 70 (Hungria -> Barein): updated case totals
$ws.Cells.Item(70, 1).Value = "Barein"
$ws.Cells.Item(70, 2).Value = 569
$ws.Cells.Item(70, 4).Value = 337
$ws.Cells.Item(70, 5).Value = 228
$ws.Cells.Item(70, 6).Value = 3
$ws.Cells.Item(70, 8).Value = 4

# Row 73 (Letonia -> Bulgaria): updated case totals
$ws.Cells.Item(73, 1).Value = "Bulgaria"
$ws.Cells.Item(73, 2).Value = 449
$ws.Cells.Item(73, 3).Value = 27
$ws.Cells.Item(73, 4).Value = 25
$ws.Cells.Item(73, 5).Value = 414
$ws.Cells.Item(73, 6).Value = 17
$ws.Cells.Item(73, 8).Value = 10

# Row 74 (Tunez -> Letonia): updated case totals
$ws.Cells.Item(74, 1).Value = "Letonia"
$ws.Cells.Item(74, 2).Value = 446
$ws.Cells.Item(74, 4).Value = 1
$ws.Cells.Item(74, 5).Value = 445
$ws.Cells.Item(74, 6).Value = 3
$ws.Cells.Item(74, 8).Value = 0

# Row 75 (Moldavia -> Tunez): updated case totals
$ws.Cells.Item(75, 1).Value = "Tunez"
$ws.Cells.Item(75, 4).Value = 5
$ws.Cells.Item(75, 5).Value = 406
$ws.Cells.Item(75, 6).Value = 10
$ws.Cells.Item(75, 8).Value = 12

# Row 76 (Bulgaria -> Moldavia): updated case totals
$ws.Cells.Item(76, 1).Value = "Moldavia"
$ws.Cells.Item(76, 2).Value = 423
$ws.Cells.Item(76, 4).Value = 23
$ws.Cells.Item(76, 5).Value = 395
$ws.Cells.Item(76, 6).Value = 44
$ws.Cells.Item(76, 8).Value = 5

# Row 94: updated case totals
$ws.Cells.Item(94, 4).Value = 64
$ws.Cells.Item(94, 5).Value = 158

# Row 102: updated case totals
$ws.Cells.Item(102, 2).Value = 187
$ws.Cells.Item(102, 3).Value = 6
$ws.Cells.Item(102, 5).Value = 173

# Row 106: updated case totals
$ws.Cells.Item(106, 5).Value = 154
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 7

# Row 107: updated case totals
$ws.Cells.Item(107, 2).Value = 148
$ws.Cells.Item(107, 3).Value = 2
$ws.Cells.Item(107, 5).Value = 124

# Row 112 (Guadalupe -> Georgia): updated case totals
$ws.Cells.Item(112, 1).Value = "Georgia"
$ws.Cells.Item(112, 2).Value = 130
$ws.Cells.Item(112, 3).Value = 13
$ws.Cells.Item(112, 4).Value = 23
$ws.Cells.Item(112, 5).Value = 107
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 8).Value = 0

# Row 113 (Montenegro -> Guadalupe): updated case totals
$ws.Cells.Item(113, 1).Value = "Guadalupe"
$ws.Cells.Item(113, 2).Value = 125
$ws.Cells.Item(113, 4).Value = 24
$ws.Cells.Item(113, 5).Value = 95
$ws.Cells.Item(113, 6).Value = 14
$ws.Cells.Item(113, 8).Value = 6

# Row 114 (Georgia -> Montenegro): updated case totals
$ws.Cells.Item(114, 1).Value = "Montenegro"
$ws.Cells.Item(114, 2).Value = 123
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 121
$ws.Cells.Item(114, 6).Value = 4
$ws.Cells.Item(114, 8).Value = 2
